# Update "想去人数" (want-to-go count) figures in column F across sheets.
# This mirrors a regenerated data refresh (gh-pages output at 456a3b4):
# the same events live on "展览"/"演出"/"本地生活" and are mirrored again
# on "全部类型", so each underlying number is bumped in both places.

$wb = $excel.ActiveWorkbook

$exhibition = $wb.Worksheets.Item("展览")
$exhibition.Cells.Item(2, 6).Value = 12697
$exhibition.Cells.Item(3, 6).Value = 7107
$exhibition.Cells.Item(25, 6).Value = 5213
$exhibition.Cells.Item(29, 6).Value = 1299
$exhibition.Cells.Item(30, 6).Value = 1299
$exhibition.Cells.Item(31, 6).Value = 47
$exhibition.Cells.Item(36, 6).Value = 589
$exhibition.Cells.Item(37, 6).Value = 64

$performance = $wb.Worksheets.Item("演出")
$performance.Cells.Item(8, 6).Value = 44

$localLife = $wb.Worksheets.Item("本地生活")
$localLife.Cells.Item(4, 6).Value = 1985

$allTypes = $wb.Worksheets.Item("全部类型")
$allTypes.Cells.Item(4, 6).Value = 1985
$allTypes.Cells.Item(6, 6).Value = 12697
$allTypes.Cells.Item(7, 6).Value = 7107
$allTypes.Cells.Item(30, 6).Value = 5213
$allTypes.Cells.Item(37, 6).Value = 1299
$allTypes.Cells.Item(38, 6).Value = 1299
$allTypes.Cells.Item(41, 6).Value = 589
$allTypes.Cells.Item(47, 6).Value = 64
